# sodan.xlsx — daily update bot: append one more day (5/17) of consultation
# figures to the "相談件数" sheet, pushing the trailing footnote row down by
# one row, and refresh the print area / view state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")
$win = $excel.ActiveWindow

# ---------------------------------------------------------------------
# 1) Insert a new row at 113. This shifts the old row 113 (the "※4/8..."
#    footnote row) down to row 114 intact, and the freshly inserted row 113
#    inherits the number formats/styles from row 112 (the row above it).
# ---------------------------------------------------------------------
$ws.Rows.Item(113).Insert()

# ---------------------------------------------------------------------
# 2) Fill the new row 113 with the day's figures (date 2020-05-17 / 43968).
# ---------------------------------------------------------------------
$ws.Range("A113").Value = 43968
$ws.Range("B113").Value = 197
$ws.Range("C113").Value = 37678
$ws.Range("D113").Value = 0
$ws.Range("E113").Value = 7584

# ---------------------------------------------------------------------
# 3) The sheet grew by one row, so the printed area grows to E115 (it was
#    already one row taller than the data, before and after).
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "=相談件数!`$A`$1:`$E`$115"
    }
}

# ---------------------------------------------------------------------
# 4) View state: zoom out a bit and move the selection to the new last
#    data row.
# ---------------------------------------------------------------------
$win.Zoom = 70
[void]$ws.Range("A113").Select()
